$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: C5 -> "Done"
$ws.Range("C5").Value = "Done"

# Row 10: C10 -> "Done", D10 -> new text about Sync Mode Config
$ws.Range("C10").Value = "Done"
$ws.Range("D10").Value = "如何Sync Mode Config？ 由很多寄存器组成"
# The new, shorter text no longer needs the extra wrapped height; let the row
# shrink back to the sheet's default row height (matches the diff dropping ht="27").
$ws.Rows.Item(10).AutoFit()

# Row 21: C21 -> "Done"
$ws.Range("C21").Value = "Done"

# Row 23: C23 -> "Done"
$ws.Range("C23").Value = "Done"

# Row 26: C26 -> "Done"
$ws.Range("C26").Value = "Done"

# Row 27: C27 -> "Done"
$ws.Range("C27").Value = "Done"

# Row 28: C28 -> "DOne" (typo as in source)
$ws.Range("C28").Value = "DOne"

# Move selection cursor to D29 (cosmetic, matches author's final cursor position)
$ws.Range("D29").Select()
